$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83-113 down to 84-114.
$ws.Rows.Item(83).Insert()

# Populate the new row 83 with the new record's data.
$ws.Range("A83").Value = 11
$ws.Range("B83").Value = "Vega Monumental Concepción"
$ws.Range("C83").Value = "Bíobío"
$ws.Range("D83").Value = 44726
$ws.Range("E83").Value = 8
$ws.Range("F83").Value = 100112021
$ws.Range("G83").Value = "Ají"
$ws.Range("H83").Value = "Inferno"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 280
$ws.Range("K83").Value = 12000
$ws.Range("L83").Value = 13000
$ws.Range("M83").Value = 12536
$ws.Range("N83").Value = "$/caja 12 kilos"
$ws.Range("O83").Value = "Región de Arica y Parinacota"
$ws.Range("P83").Value = 1045
$ws.Range("Q83").Value = 12
$ws.Range("R83").Value = "Hortaliza"
